# Apply match results for row 58 (CSK vs KKR) and let the
# dependent Rank/Total formulas recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E58").Value = 0
$ws.Range("H58").Value = 20
$ws.Range("K58").Value = 80
$ws.Range("N58").Value = 100
$ws.Range("Q58").Value = 40
$ws.Range("T58").Value = 60

$excel.Calculate()
